$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source ("Fruta, Feria Lagunitas de Puerto Montt - Frutilla") is refreshed
# weekly: two brand-new observation rows land at the top of the Frutilla
# block (row 172), and the block's prior last two rows are also re-appended
# at the very end of the sheet.

# Insert two new rows at the top of the Frutilla data block (row 172),
# pushing all existing rows in the block down by 2.
$ws.Rows.Item(172).Insert()
$ws.Rows.Item(172).Insert()

# Row 172: new weekly observation
$ws.Cells.Item(172,1).Value = 4
$ws.Cells.Item(172,2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(172,3).Value = 'Los Lagos'
$ws.Cells.Item(172,4).Value = 44642
$ws.Cells.Item(172,5).Value = 10
$ws.Cells.Item(172,6).Value = 'Fruta'
$ws.Cells.Item(172,7).Value = 100101
$ws.Cells.Item(172,8).Value = 'Berries'
$ws.Cells.Item(172,9).Value = 100112025
$ws.Cells.Item(172,10).Value = 'Frutilla'
$ws.Cells.Item(172,11).Value = 'Sin especificar'
$ws.Cells.Item(172,12).Value = 'Primera'
$ws.Cells.Item(172,13).Value = 400
$ws.Cells.Item(172,14).Value = 8000
$ws.Cells.Item(172,15).Value = 8500
$ws.Cells.Item(172,16).Value = 8250
$ws.Cells.Item(172,17).Value = '$/caja 7 kilos'
$ws.Cells.Item(172,18).Value = 'Región de La Araucanía'
$ws.Cells.Item(172,19).Value = 1179
$ws.Cells.Item(172,20).Value = 7

# Row 173: new weekly observation
$ws.Cells.Item(173,1).Value = 4
$ws.Cells.Item(173,2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(173,3).Value = 'Los Lagos'
$ws.Cells.Item(173,4).Value = 44642
$ws.Cells.Item(173,5).Value = 10
$ws.Cells.Item(173,6).Value = 'Fruta'
$ws.Cells.Item(173,7).Value = 100101
$ws.Cells.Item(173,8).Value = 'Berries'
$ws.Cells.Item(173,9).Value = 100112025
$ws.Cells.Item(173,10).Value = 'Frutilla'
$ws.Cells.Item(173,11).Value = 'Sin especificar'
$ws.Cells.Item(173,12).Value = 'Segunda'
$ws.Cells.Item(173,13).Value = 200
$ws.Cells.Item(173,14).Value = 5000
$ws.Cells.Item(173,15).Value = 5000
$ws.Cells.Item(173,16).Value = 5000
$ws.Cells.Item(173,17).Value = '$/caja 7 kilos'
$ws.Cells.Item(173,18).Value = 'Región de La Araucanía'
$ws.Cells.Item(173,19).Value = 714
$ws.Cells.Item(173,20).Value = 7

# Append two more rows at the end of the sheet (215, 216), duplicating the
# entries that were previously the last two rows of the block (213, 214).
# Row 215: new weekly observation
$ws.Cells.Item(215,1).Value = 4
$ws.Cells.Item(215,2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(215,3).Value = 'Los Lagos'
$ws.Cells.Item(215,4).Value = 44544
$ws.Cells.Item(215,5).Value = 10
$ws.Cells.Item(215,6).Value = 'Fruta'
$ws.Cells.Item(215,7).Value = 100101
$ws.Cells.Item(215,8).Value = 'Berries'
$ws.Cells.Item(215,9).Value = 100112025
$ws.Cells.Item(215,10).Value = 'Frutilla'
$ws.Cells.Item(215,11).Value = 'Sin especificar'
$ws.Cells.Item(215,12).Value = 'Primera'
$ws.Cells.Item(215,13).Value = 800
$ws.Cells.Item(215,14).Value = 10000
$ws.Cells.Item(215,15).Value = 11000
$ws.Cells.Item(215,16).Value = 10500
$ws.Cells.Item(215,17).Value = '$/bandeja 7 kilos'
$ws.Cells.Item(215,18).Value = 'Provincia de Melipilla'
$ws.Cells.Item(215,19).Value = 1500
$ws.Cells.Item(215,20).Value = 7

# Row 216: new weekly observation
$ws.Cells.Item(216,1).Value = 4
$ws.Cells.Item(216,2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(216,3).Value = 'Los Lagos'
$ws.Cells.Item(216,4).Value = 44544
$ws.Cells.Item(216,5).Value = 10
$ws.Cells.Item(216,6).Value = 'Fruta'
$ws.Cells.Item(216,7).Value = 100101
$ws.Cells.Item(216,8).Value = 'Berries'
$ws.Cells.Item(216,9).Value = 100112025
$ws.Cells.Item(216,10).Value = 'Frutilla'
$ws.Cells.Item(216,11).Value = 'Sin especificar'
$ws.Cells.Item(216,12).Value = 'Primera'
$ws.Cells.Item(216,13).Value = 800
$ws.Cells.Item(216,14).Value = 8500
$ws.Cells.Item(216,15).Value = 9000
$ws.Cells.Item(216,16).Value = 8750
$ws.Cells.Item(216,17).Value = '$/caja 7 kilos'
$ws.Cells.Item(216,18).Value = 'Región de La Araucanía'
$ws.Cells.Item(216,19).Value = 1250
$ws.Cells.Item(216,20).Value = 7
